$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" — zero out a handful of category totals for row
# 4, 14, 18 and 21 (their values are being moved out of this month's sheet),
# and correct the "x de 27" progress labels in row 29 to "0 de 27".
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("K4").Value = 0
$wsGrupo.Range("K14").Value = 0
$wsGrupo.Range("L14").Value = 0
$wsGrupo.Range("L18").Value = 0
$wsGrupo.Range("C21").Value = 0
$wsGrupo.Range("M21").Value = 0
$wsGrupo.Range("N21").Value = 0

$wsGrupo.Range("C29").Value = "0 de 27"
$wsGrupo.Range("K29").Value = "0 de 27"
$wsGrupo.Range("L29").Value = "0 de 27"
$wsGrupo.Range("M29").Value = "0 de 27"
$wsGrupo.Range("N29").Value = "0 de 27"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL" — monthly rollover: febrero/marzo/abril/mayo becomes
# marzo/abril/mayo/junio, so every month column of data shifts one column to
# the left (C<-D, D<-E, E<-F, F<-0). Column widths shift the same way.
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("C1").Value = "marzo"
$wsMensual.Range("D1").Value = "abril"
$wsMensual.Range("E1").Value = "mayo"
$wsMensual.Range("F1").Value = "junio"

# Note: ColumnWidth -> stored sheet width has a constant +5/6 offset in
# this engine's MDW-7 pixel-grid conversion, so subtract it here to land
# exactly on the target stored widths (11, 13, 11).
$wsMensual.Columns.Item(4).ColumnWidth = 11 - 5/6
$wsMensual.Columns.Item(5).ColumnWidth = 13 - 5/6
$wsMensual.Columns.Item(6).ColumnWidth = 11 - 5/6

$wsMensual.Range("C3").Value = 0

$wsMensual.Range("E4").Value = 419.13
$wsMensual.Range("F4").Value = 0

$wsMensual.Range("C14").Value = 1444.13
$wsMensual.Range("D14").Value = 0
$wsMensual.Range("E14").Value = 3122.02
$wsMensual.Range("F14").Value = 0

$wsMensual.Range("D16").Value = 226.8
$wsMensual.Range("E16").Value = 0

$wsMensual.Range("E18").Value = 6725.74
$wsMensual.Range("F18").Value = 0

$wsMensual.Range("C19").Value = 800.79
$wsMensual.Range("D19").Value = 1126
$wsMensual.Range("E19").Value = 0

$wsMensual.Range("C21").Value = 738.55
$wsMensual.Range("D21").Value = 0
$wsMensual.Range("E21").Value = 1994.73
$wsMensual.Range("F21").Value = 0

$wsMensual.Range("C27").Value = 0

$wsMensual.Range("C29").Value = 2983.47
$wsMensual.Range("D29").Value = 1352.8
$wsMensual.Range("E29").Value = 12261.62
$wsMensual.Range("F29").Value = 0
